$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 807.38464
$ws.Range("I103").Value = 1163.3334
$ws.Range("J103").Value = 502.2857
$ws.Range("K103").Value = 3490.0002
$ws.Range("L103").Value = 1506.8571
$ws.Range("M103").Value = -2904.0002
$ws.Range("N103").Value = -2678.8571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2840.6226
$ws.Range("I132").Value = 2873.319
$ws.Range("J132").Value = 2584.5
$ws.Range("K132").Value = 8619.957
$ws.Range("L132").Value = 7753.5
$ws.Range("M132").Value = -6089.957
$ws.Range("N132").Value = -12813.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 55556276
$ws.Range("I135").Value = 26316344
$ws.Range("J135").Value = 125001120
$ws.Range("K135").Value = 236847096
$ws.Range("L135").Value = 1125010080
$ws.Range("M135").Value = -236844561
$ws.Range("N135").Value = -1125015150

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1262.5
$ws.Range("J137").Value = 1446
$ws.Range("L137").Value = 4338
$ws.Range("N137").Value = -9438

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 47622456
$ws.Range("I74").Value = 52635252
$ws.Range("K74").Value = 52635252
$ws.Range("M74").Value = -52634378

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 47622456
$ws.Range("I77").Value = 52635252
$ws.Range("K77").Value = 263176260
$ws.Range("M77").Value = -263171892

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2942515
$ws.Range("I132").Value = 2942515
$ws.Range("K132").Value = 8827545
$ws.Range("M132").Value = -8825015

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 50000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("M123").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 49999.125
$ws.Range("J126").Value = 49999.125
$ws.Range("L126").Value = 49999.125
$ws.Range("N126").Value = -59879.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 16452747
$ws.Range("I134").Value = 17587206
$ws.Range("J134").Value = 3095
$ws.Range("K134").Value = 52761618
$ws.Range("L134").Value = 9285
$ws.Range("M134").Value = -52759083
$ws.Range("N134").Value = -14355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3771.9707
$ws.Range("I31").Value = 5141.722
$ws.Range("J31").Value = 2231
$ws.Range("K31").Value = 5141.722
$ws.Range("L31").Value = 2231
$ws.Range("M31").Value = -4846.722
$ws.Range("N31").Value = -2821

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3771.9707
$ws.Range("I34").Value = 5141.722
$ws.Range("J34").Value = 2231
$ws.Range("K34").Value = 5141.722
$ws.Range("L34").Value = 2231
$ws.Range("M34").Value = -4939.722
$ws.Range("N34").Value = -2635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 44998.332
$ws.Range("J50").Value = 44998.332
$ws.Range("L50").Value = 44998.332
$ws.Range("N50").Value = -46248.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 25631.818
$ws.Range("I60").Value = 8490
$ws.Range("K60").Value = 8490
$ws.Range("M60").Value = -7979

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 69500
$ws.Range("J98").Value = 69500
$ws.Range("L98").Value = 69500
$ws.Range("N98").Value = -73992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2673.7083
$ws.Range("I122").Value = 3012.8667
$ws.Range("K122").Value = 9038.6001
$ws.Range("M122").Value = -6588.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H123").Value = 80000
$ws.Range("J123").Value = 80000
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -89800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 45163
$ws.Range("J125").Value = 45163
$ws.Range("L125").Value = 45163
$ws.Range("N125").Value = -50083

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 50000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 50000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920
$ws.Range("M127").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 41671270
$ws.Range("I132").Value = 50004624
$ws.Range("J132").Value = 4510.25
$ws.Range("K132").Value = 150013872
$ws.Range("L132").Value = 13530.75
$ws.Range("M132").Value = -150011342
$ws.Range("N132").Value = -18590.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9618567
$ws.Range("I134").Value = 10002910
$ws.Range("K134").Value = 30008730
$ws.Range("M134").Value = -30006195

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 321
$ws.Range("I14").Value = 321
$ws.Range("K14").Value = 963
$ws.Range("M14").Value = -790

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 81.40000000000001
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 101.15385
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 101.15385
$ws.Range("M2").Value = 53
$ws.Range("N2").Value = -327.15385

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 58845.777
$ws.Range("I113").Value = 74613.71000000001
$ws.Range("K113").Value = 74613.71000000001
$ws.Range("M113").Value = -72443.71000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2726.4443
$ws.Range("I126").Value = 2842.25
$ws.Range("K126").Value = 8526.75
$ws.Range("M126").Value = -6056.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 150
$ws.Range("I10").Value = 150
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 150
$ws.Range("L10").Value = 150
$ws.Range("M10").Value = -10
$ws.Range("N10").Value = -430

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 728.1111
$ws.Range("I55").Value = 260.5
$ws.Range("K55").Value = 260.5
$ws.Range("M55").Value = -87.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3852.4707
$ws.Range("I136").Value = 3852.4707
$ws.Range("K136").Value = 11557.4121
$ws.Range("M136").Value = -9007.4121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5939.8
$ws.Range("I62").Value = 2851.25
$ws.Range("K62").Value = 2851.25
$ws.Range("M62").Value = -2227.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5939.8
$ws.Range("I65").Value = 2851.25
$ws.Range("K65").Value = 14256.25
$ws.Range("M65").Value = -11136.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2679.8
$ws.Range("I126").Value = 3125
$ws.Range("K126").Value = 9375
$ws.Range("M126").Value = -6905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 45466690
$ws.Range("I132").Value = 62506760
$ws.Range("K132").Value = 187520280
$ws.Range("M132").Value = -187517750

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 50002740
$ws.Range("I136").Value = 55557490
$ws.Range("K136").Value = 166672470
$ws.Range("M136").Value = -166669920
